# Apply "weekly" row re-shuffle to the Perejil (Hortaliza) data sheet.
# The columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado), N (Unidad de comercializacion),
# P (Precio $/Kg) and Q (Kg o Unidades) for data rows 2-29 are permuted
# (each destination row receives the full tuple of values from a specific
# source row), while all other columns (A,B,C,E,F,G,H,I,O,R) and row 1
# (headers) stay untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that move together as a unit for every data row.
$cols = @("D", "J", "K", "L", "M", "N", "P", "Q")

# Mapping: destination row -> source row (values are read from the
# ORIGINAL/before state of the worksheet before any writes happen).
$rowMap = @{
    2  = 29
    3  = 27
    4  = 19
    5  = 13
    6  = 20
    7  = 22
    8  = 7
    9  = 9
    10 = 4
    11 = 25
    12 = 12
    13 = 21
    14 = 23
    15 = 16
    16 = 14
    17 = 3
    18 = 26
    19 = 6
    20 = 24
    21 = 28
    22 = 8
    23 = 17
    24 = 2
    25 = 10
    26 = 11
    27 = 15
    28 = 18
    29 = 5
}

# 1) Snapshot the original values for every column/row we might need,
#    BEFORE writing anything, so overwrites don't clobber source data.
#    (Value2 is used instead of Value: it round-trips cleanly through
#    PS variables/hashtables in this COM-interop host.)
$snapshot = @{}
foreach ($row in 2..29) {
    $rowData = @{}
    foreach ($col in $cols) {
        $rowData[$col] = $ws.Range("$col$row").Value2
    }
    $snapshot[$row] = $rowData
}

# 2) Write the permuted values back into the sheet.
foreach ($destRow in 2..29) {
    $srcRow = $rowMap[$destRow]
    $srcData = $snapshot[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$destRow").Value2 = $srcData[$col]
    }
}
